$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7 (this shifts existing rows 7:100 down to 8:101,
# carrying their values/formats with them, and extends the used range to
# A1:T101).
$ws.Rows("7:7").Insert()

# Populate the newly inserted (blank) row 7 with a fresh data record. The
# static descriptive columns (A,B,C,E,F,G,H,I,J,K,L,M,Q,R,T) mirror the
# record that used to sit in row 7 (now in row 8); only the date and the
# price columns (D,N,O,P,S) carry new values.
$ws.Range("A7").Value = 11
$ws.Range("B7").Value = "Vega Monumental Concepción"
$ws.Range("C7").Value = "Bíobío"
$ws.Range("D7").Value = "2022-02-24"
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100108
$ws.Range("H7").Value = "Tropicales y subtropicales"
$ws.Range("I7").Value = 100108002
$ws.Range("J7").Value = "Mango"
$ws.Range("K7").Value = "Sin especificar"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 200
$ws.Range("N7").Value = 7000
$ws.Range("O7").Value = 8000
$ws.Range("P7").Value = 7500
$ws.Range("Q7").Value = "`$/bandeja 4 kilos"
$ws.Range("R7").Value = "Perú"
$ws.Range("S7").Value = 1875
$ws.Range("T7").Value = 4
